$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of rows 7 through 15 (columns A:W), keeping formatting/styles.
$ws.Range("A7:W15").ClearContents()

# Update the view: scroll so column B is the left-most visible column,
# and select D4.
$ws.Range("D4").Select()
$ws.Application.ActiveWindow.ScrollColumn = 2
